{"js": "// Fix: the paragraph holding the \"{{loop_falla}}\" merge field loop tag was\n// rendering with no explicit font, which made the report repeat images\n// when the code was the same. Force the run (and the paragraph mark) to\n// Arial 11pt so every iteration renders consistently.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.indexOf(\"{{loop_falla}}\") !== -1) {\n    const font = paragraph.font;\n    font.name = \"Arial\";\n    font.size = 11;\n    // Keep the complex-script size attribute (w:szCs) in sync too.\n    font.sizeBidirectional = 11;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix: the paragraph holding the \"{{loop_falla}}\" merge field loop tag was\n# rendering with no explicit font, which made the report repeat images\n# when the code was the same. Force the run (and the paragraph mark) to\n# Arial 11pt so every iteration renders consistently.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $rng = $p.Range\n    if ($rng.Text -like \"*{{loop_falla}}*\") {\n        $rng.Font.Name = \"Arial\"\n        $rng.Font.Size = 11\n        # Keep the complex-script size attribute (w:szCs) in sync too.\n        $rng.Font.SizeBi = 11\n    }\n}\n"}
